$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.037.85'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.829.89'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.13'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6232'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.60%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07537'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.56'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +6.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2911'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.79'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07637'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.828.68'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6646'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.34'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009102'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +8.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.000'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '28.808.04'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.083.37'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '224.58'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.98%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.192'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.000'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.50'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.386'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1355'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.83'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.497'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.048'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.41%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.029'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.17%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.199'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05200'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.833'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.152'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7316'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.610'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.289.91'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.754'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.376'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +6.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8916'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.44'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.981.21'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.45'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3969'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.74%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.852'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.644'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -5.78%  '

Write-Host "Applied all changes"